$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date-range banner) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Helper functions to flip a cell between the "0"/"***.*" text placeholder
#     style (s=14, shared text) and a plain numeric style, while reusing an
#     existing same-column donor cell (one untouched by this edit) so the
#     resulting cell style id matches the target exactly. ---
function Set-PlaceholderText($ref, $donorRef, $text) {
    $dst = $ws.Range($ref)
    $dst.NumberFormat = "@"
    $dst.Value = $text
    $src = $ws.Range($donorRef)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

function Set-NumberFromText($ref, $donorRef, $num) {
    $dst = $ws.Range($ref)
    $dst.Value = $num
    $src = $ws.Range($donorRef)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# --- Cells that flip from a number to the "0" / "***.*" text placeholder ---
Set-PlaceholderText "C14" "C23" "0"
Set-PlaceholderText "D18" "D23" "0"
Set-PlaceholderText "E18" "E23" "***.*"
Set-PlaceholderText "C20" "C23" "0"
Set-PlaceholderText "D20" "D23" "0"
Set-PlaceholderText "E20" "E23" "***.*"
Set-PlaceholderText "F28" "F23" "0"
Set-PlaceholderText "F29" "F23" "0"

# --- Cells that flip from the text placeholder to a real number ---
Set-NumberFromText "D15" "D17" 1
Set-NumberFromText "E15" "E17" -100
Set-NumberFromText "D26" "D17" 2
Set-NumberFromText "E26" "E17" -100
Set-NumberFromText "C30" "C17" 1
Set-NumberFromText "F30" "F17" 1

$excel.CutCopyMode = $false

# --- Plain numeric updates (style unchanged) ---
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = -57.142857142857
$ws.Range("M15").Value = -40
$ws.Range("N15").Value = -62.5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 126
$ws.Range("J16").Value = 152
$ws.Range("K16").Value = -17.105263157894
$ws.Range("L16").Value = 0.8
$ws.Range("M16").Value = -5.263157894736
$ws.Range("N16").Value = -90.374331550802
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 169
$ws.Range("J17").Value = 175
$ws.Range("K17").Value = -3.428571428571
$ws.Range("L17").Value = 19.014084507042
$ws.Range("M17").Value = 37.398373983739
$ws.Range("N17").Value = -63.888888888888
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 127
$ws.Range("K18").Value = -36.180904522613
$ws.Range("L18").Value = -33.157894736842
$ws.Range("M18").Value = -28.651685393258
$ws.Range("N18").Value = -93.753074274471
$ws.Range("C19").Value = 55
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = 27.906976744186
$ws.Range("F19").Value = 185
$ws.Range("G19").Value = 177
$ws.Range("H19").Value = 4.519774011299
$ws.Range("I19").Value = 1759
$ws.Range("J19").Value = 1890
$ws.Range("K19").Value = -6.931216931216
$ws.Range("L19").Value = 44.298605414274
$ws.Range("M19").Value = 5.836341756919
$ws.Range("N19").Value = -75.319208643187
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 77
$ws.Range("K20").Value = -31.25
$ws.Range("L20").Value = 8.450704225352
$ws.Range("M20").Value = 83.333333333333
$ws.Range("N20").Value = -83.476394849785
$ws.Range("C21").Value = 63
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = 36.956521739130
$ws.Range("F21").Value = 219
$ws.Range("G21").Value = 210
$ws.Range("H21").Value = 4.285714285714
$ws.Range("I21").Value = 2270
$ws.Range("J21").Value = 2550
$ws.Range("K21").Value = -10.980392156862
$ws.Range("L21").Value = 28.393665158371
$ws.Range("M21").Value = 5.287569573283
$ws.Range("N21").Value = -80.153873054729
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -11.111111111111
$ws.Range("I22").Value = 76
$ws.Range("J22").Value = 68
$ws.Range("K22").Value = 11.764705882352
$ws.Range("L22").Value = 61.702127659574
$ws.Range("M22").Value = 22.580645161290
$ws.Range("C24").Value = 66
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = 15.789473684210
$ws.Range("F24").Value = 215
$ws.Range("G24").Value = 259
$ws.Range("H24").Value = -16.988416988417
$ws.Range("I24").Value = 2597
$ws.Range("J24").Value = 2631
$ws.Range("K24").Value = -1.292284302546
$ws.Range("L24").Value = 38.211814795103
$ws.Range("M24").Value = 40.454299621417
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -8.333333333333
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 9.756097560975
$ws.Range("I25").Value = 640
$ws.Range("J25").Value = 541
$ws.Range("K25").Value = 18.299445471349
$ws.Range("L25").Value = 48.491879350348
$ws.Range("M25").Value = 47.465437788018
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = -37.142857142857
$ws.Range("L26").Value = -29.032258064516
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 98
$ws.Range("J27").Value = 98
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 19.512195121951
$ws.Range("N28").Value = -61.538461538461
$ws.Range("N29").Value = -60
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 20
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = 17.647058823529
